$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to remain text so numeric-looking strings (e.g. "0.999", "2.60")
# keep their exact textual representation instead of being coerced to floats.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "63.364.09"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "3.341.33"
$ws.Range("E3").Value = "  +1.32%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "529.46"
$ws.Range("E5").Value = "  +1.60%  "
$ws.Range("D6").Value = "172.67"
$ws.Range("E6").Value = "  -5.53%  "
$ws.Range("D7").Value = "0.595"
$ws.Range("E7").Value = "  -0.39%  "
$ws.Range("D8").Value = "3.333.94"
$ws.Range("E8").Value = "  +0.90%  "
$ws.Range("E9").Value = "  -0.03%  "
$ws.Range("D10").Value = "0.608"
$ws.Range("E10").Value = "  -2.02%  "
$ws.Range("D11").Value = "53.13"
$ws.Range("E11").Value = "  -8.85%  "
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("D13").Value = "0.0000255"
$ws.Range("E13").Value = "  -0.40%  "
$ws.Range("D14").Value = "9.13"
$ws.Range("E14").Value = "  -0.12%  "
$ws.Range("D15").Value = "3.882.50"
$ws.Range("E15").Value = "  +1.71%  "
$ws.Range("D16").Value = "3.343.55"
$ws.Range("E16").Value = "  +1.47%  "
$ws.Range("E17").Value = "  -1.13%  "
$ws.Range("D18").Value = "17.45"
$ws.Range("E18").Value = "  -1.69%  "
$ws.Range("D19").Value = "63.282.15"
$ws.Range("E19").Value = "  -0.51%  "
$ws.Range("E20").Value = "  +1.95%  "
$ws.Range("D21").Value = "0.964"
$ws.Range("E21").Value = "  +1.31%  "
$ws.Range("D22").Value = "371.83"
$ws.Range("E22").Value = "  -0.04%  "
$ws.Range("D23").Value = "4.10"
$ws.Range("E23").Value = "  +6.70%  "
$ws.Range("D24").Value = "11.32"
$ws.Range("E24").Value = "  +0.44%  "
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "6.16"
$ws.Range("E27").Value = "  +2.91%  "
$ws.Range("E28").Value = "  +1.80%  "
$ws.Range("D29").Value = "11.29"
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("D30").Value = "8.28"
$ws.Range("E30").Value = "  -0.67%  "
$ws.Range("D31").Value = "28.86"
$ws.Range("E31").Value = "  +0.87%  "
$ws.Range("D32").Value = "635.77"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").Value = "6.42"
$ws.Range("E33").Value = "  -4.08%  "
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").Value = "57.76"
$ws.Range("E36").Value = "  -2.84%  "
$ws.Range("D38").Value = "36.83"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").Value = "0.380"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "0.0₃0723"
$ws.Range("E40").Value = "  +10.78%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").Value = "2.61"
$ws.Range("E42").Value = "  +7.71%  "
$ws.Range("D43").Value = "0.124"
$ws.Range("E43").Value = "  -1.00%  "
$ws.Range("D44").Value = "2.915.11"
$ws.Range("E44").Value = "  -2.63%  "
$ws.Range("E45").Value = "  +5.99%  "
$ws.Range("E46").Value = "  +3.04%  "
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("D48").Value = "2.60"
$ws.Range("E48").Value = "  -2.93%  "
$ws.Range("E49").Value = "  +4.00%  "
$ws.Range("E50").Value = "  -0.43%  "
$ws.Range("D51").Value = "135.66"
$ws.Range("E51").Value = "  +3.52%  "
